# Update pretestSummary.xlsx template:
#  - repoint jx:each loop from summary.preTestData/course to summary.uploadRecords/r
#  - replace the pre-test start/end/duration columns with file-upload columns
#    (FILE PATH / FILE NAME / UPLOAD DATE) and adjust headers/values
#  - widen columns B and D to fit the new content
#  - update selection + the A8 comment describing the loop

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 7) ---
# NOTE: these values use literal ${...} template placeholders, which must be
# single-quoted so PowerShell does not try to interpolate them as variables.
$ws.Range("A7").Value = 'USER ID'
$ws.Range("B7").Value = 'FILE PATH'
$ws.Range("C7").Value = 'FILE NAME'
$ws.Range("D7").Value = 'UPLOAD DATE'

# --- Data row (row 8), driven by the new loop variable "r" ---
$ws.Range("A8").Value = '${r.userId}'
$ws.Range("B8").Value = '${r.filePath}'
$ws.Range("C8").Value = '${r.fileName}'
$ws.Range("D8").Value = '${utils:dateFmt(r.uploadDate, "yyyy-MM-dd HH:mm:ss")}'

# --- Column widths: B -> 35, D -> 25.75 (C is left untouched at 23.25) ---
# (ColumnWidth adds ~5/7 of a character internally, so compensate on input)
$ws.Columns.Item(2).ColumnWidth = 35 - 5/7
$ws.Columns.Item(4).ColumnWidth = 25.75 - 5/7

# --- Selection moves from A6 to D8 ---
$ws.Range("D8").Select()

# --- Update the A8 comment text to describe the new loop ---
$cmt = $ws.Range("A8").Comment
$cmt.Text('jx:each(items="summary.uploadRecords" var="r" lastCell="D8")')
